# feat: add 2022-Q4 data
#
# - Insert a new sheet "2022-Q4" right before "2022-Q3", populated with the
#   Q4 fund-holdings table (same layout as the existing quarterly sheets).
# - Insert a new leading data row in "总计" summarizing the Q4 quarter,
#   pushing the existing 2022-Q3 / 2022-Q2 summary rows down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the "2022-Q4" worksheet by cloning "2022-Q3" (same header row /
#    column-A style already in place) and placing the clone right before
#    it, then overwrite its data with the Q4 figures.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
# After Copy(), $q3 rebinds to the freshly-inserted clone (it now occupies
# the position the selector resolves to); the original "2022-Q3" sheet got
# pushed one slot to the right and is looked up again below.
$q4 = $q3
$q4.Name = "2022-Q4"

# The clone only has 2 data rows (rows 2:3); extend the same per-cell
# formatting (styled col A, plain B:H) down to rows 4:6 for the 3 extra
# funds in the Q4 table.
$q4.Range("A3:H3").Copy()
$q4.Range("A4:H6").PasteSpecial(-4122)

$q4Rows = @(
    @(0, "014135", "中欧金安量化混合A",               "7.09", "90.01", "0.78", "0.0553", 2),
    @(1, "001990", "中欧数据挖掘多因子灵活配置混合A", "3.26", "90.14", "0.82", "0.0267", 2),
    @(2, "004234", "中欧数据挖掘多因子灵活配置混合C", "1.93", "90.14", "0.82", "0.0158", 2),
    @(3, "014136", "中欧金安量化混合C",               "1.07", "90.01", "0.78", "0.0083", 2),
    @(4, "005167", "嘉实润泽量化一年定期开放混合",     "0.56", "27.25", "0.67", "0.0038", 1)
)

foreach ($row in $q4Rows) {
    $r = [int]$row[0] + 2

    $q4.Cells.Item($r, 1).Value = [int]$row[0]

    # Fund code / name / size / position / ratio / market-value are all
    # stored as plain text (leading zeros in fund codes, etc.) - force
    # text with a leading apostrophe, then strip the quote-prefix style
    # it introduces so the cell keeps the (unstyled) look of its peers.
    for ($col = 2; $col -le 7; $col++) {
        $cell = $q4.Cells.Item($r, $col)
        $cell.Value = "'" + $row[$col - 1]
        $cell.Style = "Normal"
    }

    $q4.Cells.Item($r, 8).Value = [int]$row[7]
}

# ---------------------------------------------------------------------
# 2) Insert a new leading row in "总计" for the 2022-Q4 summary, shifting
#    the existing rows (2022-Q3, 2022-Q2) down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend the existing row-3 formatting down to row 4 (new row), so row 4
# gets the same per-cell styles as rows 2/3 (col A styled, B:D plain).
$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4122)

# Shift 2022-Q2 content down into row 4, 2022-Q3 into row 3 (re-asserted,
# unchanged values), and write the new 2022-Q4 summary into row 2.
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q2"
$total.Cells.Item(4, 3).Value = 1
$total.Cells.Item(4, 4).Value = 0

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 2
$total.Cells.Item(3, 4).Value = 0

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.11
